# Corrects the JF code values for utilities (row 10) in Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C10").Value  = 9.603478434330045
$ws.Range("I10").Value  = 2.2017940289978917
$ws.Range("K10").Value  = 0.0
$ws.Range("P10").Value  = 4.572825479844843
$ws.Range("R10").Value  = 2.843261580753183
$ws.Range("S10").Value  = 2.080435302990134
$ws.Range("T10").Value  = 2.395967990610304
$ws.Range("U10").Value  = 2.066837232407834
$ws.Range("V10").Value  = 0.0
$ws.Range("W10").Value  = 2.1671201072813897
$ws.Range("X10").Value  = 4.137128678467276
$ws.Range("Y10").Value  = 6.202875002252188
$ws.Range("Z10").Value  = 2.451446265356708
$ws.Range("AA10").Value = 0.0
$ws.Range("AI10").Value = 5.719589512720284
$ws.Range("AL10").Value = 3.0409029345372454
$ws.Range("AN10").Value = 2.1983266368262413
$ws.Range("AO10").Value = 2.4271745201551562
$ws.Range("AR10").Value = 4.591061900435251
$ws.Range("AS10").Value = 3.032329457853188
$ws.Range("AT10").Value = 1.6634860761918413
$ws.Range("AU10").Value = 2.128978793393237
$ws.Range("AV10").Value = 0.0
$ws.Range("AY10").Value = 4.512139347405864
$ws.Range("BA10").Value = 8.481848004303611
